# Generate Report for Handoff
#
# The "b.md" row (row 3) on every sheet is updated to reflect that a new
# handoff has just been generated:
#   - Status moves from "Handed back: in sync with en-US" to "Ready for handoff"
#   - A new (newer) handoff xlf file name + timestamp is recorded
#   - An error/detail message is now populated, explaining that the handback
#     file that was previously received is stale relative to the new handoff
#
# Cell values are assigned in the same left-to-right / top-to-bottom / sheet
# order that the new text first appears in, so that new shared-string table
# entries come out in the right order.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2923d7408876db1b5c415f5157c960e91909e988/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9c964ca8fb1a82331fe6f6b54f297e7016b48d60/e2e/b.md."

# --- Overview sheet: row 3 (b.md) ---
$ws1.Range("E3").Value = "Ready for handoff"
$ws1.Range("F3").Value = "Ready for handoff"
$ws1.Range("G3").Value = "2016-10-18 12:20:45"

# --- zh-cn sheet: row 3 (b.md) ---
$ws2.Range("C3").Value = "Ready for handoff"
$ws2.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$ws2.Range("H3").Value = "2016-10-18 12:20:34"
$ws2.Range("P3").Value = $errorDetail
$ws2.Columns.Item(16).ColumnWidth = 39.17

# --- de-de sheet: row 3 (b.md) ---
$ws3.Range("C3").Value = "Ready for handoff"
$ws3.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$ws3.Range("H3").Value = "2016-10-18 12:20:45"
$ws3.Range("P3").Value = $errorDetail
$ws3.Columns.Item(16).ColumnWidth = 39.17
